$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 9: header row for the new URL lookup table (bold/centered like rows 1 & 6) ---
$ws.Range("A9").Value = "type"
$ws.Range("B9").Value = "genome"
$ws.Range("C9").Value = "URL"
$ws.Range("H9").Value = "[this line is not part of the file format itself]"

# Copy the header formatting from row 1 onto the new header cells individually
# (keeps only the cells that actually carry a value, matching the existing
# sparse-row layout used throughout this sheet).
$ws.Range("A1").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("B1").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("C1").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("H1").Copy()
$ws.Range("H9").PasteSpecial(-4122)

# --- Row 10: Barley URL entry ---
$ws.Range("A10").Value = "URL"
$ws.Range("B10").Value = "Barley"
$ws.Range("C10").Value = "http://penguin.scri.ac.uk/paul/germinate/germinate_development/app/flapjack/flapjack_search/search.pl?marker="

# --- Row 11: Rice URL entry ---
$ws.Range("A11").Value = "URL"
$ws.Range("B11").Value = "Rice"
$ws.Range("C11").Value = "http://rice.plantbiology.msu.edu/cgi-bin/gbrowse/rice/?name="

# Move / update the active selection to match the post-edit state
$ws.Range("H9").Select()
